$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 997.6
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 997
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 997
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -1347
# Row 80
$ws.Range("H80").Value = 520.8421
$ws.Range("J80").Value = 594.44446
$ws.Range("L80").Value = 1783.33338
$ws.Range("N80").Value = -3779.33338
# Row 83
$ws.Range("H83").Value = 520.8421
$ws.Range("J83").Value = 594.44446
$ws.Range("L83").Value = 5350.00014
$ws.Range("N83").Value = -15334.00014
# Row 113
$ws.Range("H113").Value = 10418766
$ws.Range("I113").Value = 3100
$ws.Range("K113").Value = 3100
$ws.Range("M113").Value = 154
# Row 129
$ws.Range("H129").Value = 1047.7656
$ws.Range("J129").Value = 1125.0962
$ws.Range("L129").Value = 3375.2886
$ws.Range("N129").Value = -13375.2886
# Row 132
$ws.Range("H132").Value = 1795.6072
$ws.Range("I132").Value = 1556.0526
$ws.Range("J132").Value = 2301.3333
$ws.Range("K132").Value = 4668.1578
$ws.Range("L132").Value = 6903.999899999999
$ws.Range("M132").Value = -2138.1578
$ws.Range("N132").Value = -11963.9999
# Row 138
$ws.Range("H138").Value = 2660.9143
$ws.Range("I138").Value = 1184.5186
$ws.Range("J138").Value = 3587.9534
$ws.Range("K138").Value = 3553.5558
$ws.Range("L138").Value = 10763.8602
$ws.Range("M138").Value = 1586.4442
$ws.Range("N138").Value = -21043.8602

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 648.4815
$ws.Range("I2").Value = 664.24
$ws.Range("J2").Value = 451.5
$ws.Range("K2").Value = 664.24
$ws.Range("L2").Value = 451.5
$ws.Range("M2").Value = -551.24
$ws.Range("N2").Value = -677.5
# Row 32
$ws.Range("H32").Value = 3583.37
$ws.Range("I32").Value = 2500.277
$ws.Range("J32").Value = 8871.412
$ws.Range("K32").Value = 2500.277
$ws.Range("L32").Value = 8871.412
$ws.Range("M32").Value = -2213.277
$ws.Range("N32").Value = -9445.412
# Row 45
$ws.Range("H45").Value = 13097
$ws.Range("I45").Value = 14825.143
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 14825.143
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -14448.143
$ws.Range("N45").Value = -1754
# Row 61
$ws.Range("H61").Value = 3971.8975
$ws.Range("I61").Value = 4551.533
$ws.Range("J61").Value = 2039.7778
$ws.Range("K61").Value = 4551.533
$ws.Range("L61").Value = 2039.7778
$ws.Range("M61").Value = -4339.533
$ws.Range("N61").Value = -2463.7778
# Row 74
$ws.Range("H74").Value = 932.1539
$ws.Range("I74").Value = 758.2857
$ws.Range("J74").Value = 1464.625
$ws.Range("K74").Value = 758.2857
$ws.Range("L74").Value = 1464.625
$ws.Range("M74").Value = 115.7143
$ws.Range("N74").Value = -3212.625
# Row 77
$ws.Range("H77").Value = 932.1539
$ws.Range("I77").Value = 758.2857
$ws.Range("J77").Value = 1464.625
$ws.Range("K77").Value = 3791.4285
$ws.Range("L77").Value = 7323.125
$ws.Range("M77").Value = 576.5715
$ws.Range("N77").Value = -16059.125
# Row 110
$ws.Range("H110").Value = 1250
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 1000
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 1045
$ws.Range("N110").Value = -5590
# Row 116
$ws.Range("H116").Value = 648.4815
$ws.Range("I116").Value = 664.24
$ws.Range("J116").Value = 451.5
$ws.Range("K116").Value = 664.24
$ws.Range("L116").Value = 451.5
$ws.Range("M116").Value = 1629.76
$ws.Range("N116").Value = -5039.5
# Row 122
$ws.Range("H122").Value = 2565107.2
$ws.Range("I122").Value = 3206045
$ws.Range("K122").Value = 9618135
$ws.Range("M122").Value = -9615685
# Row 132
$ws.Range("H132").Value = 3977.6875
$ws.Range("I132").Value = 2681.0476
$ws.Range("K132").Value = 8043.1428
$ws.Range("M132").Value = -5513.1428
# Row 136
$ws.Range("H136").Value = 3971.8975
$ws.Range("I136").Value = 4551.533
$ws.Range("J136").Value = 2039.7778
$ws.Range("K136").Value = 13654.599
$ws.Range("L136").Value = 6119.3334
$ws.Range("M136").Value = -11104.599
$ws.Range("N136").Value = -11219.3334

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 648.4815
$ws.Range("I3").Value = 664.24
$ws.Range("J3").Value = 451.5
$ws.Range("K3").Value = 664.24
$ws.Range("L3").Value = 451.5
$ws.Range("M3").Value = -550.24
$ws.Range("N3").Value = -679.5
# Row 99
$ws.Range("H99").Value = 47620236
$ws.Range("I99").Value = 58824452
$ws.Range("J99").Value = 2324.75
$ws.Range("K99").Value = 58824452
$ws.Range("L99").Value = 2324.75
$ws.Range("M99").Value = -58822954
$ws.Range("N99").Value = -5320.75
# Row 105
$ws.Range("H105").Value = 12201658
$ws.Range("I105").Value = 20842884
$ws.Range("J105").Value = 2280.5881
$ws.Range("K105").Value = 20842884
$ws.Range("L105").Value = 2280.5881
$ws.Range("M105").Value = -20841137
$ws.Range("N105").Value = -5774.5881
# Row 107
$ws.Range("H107").Value = 83334630
$ws.Range("I107").Value = 200001200
$ws.Range("J107").Value = 1362.2858
$ws.Range("K107").Value = 200001200
$ws.Range("L107").Value = 1362.2858
$ws.Range("M107").Value = -199999280
$ws.Range("N107").Value = -5202.2858
# Row 134
$ws.Range("H134").Value = 6360.8887
$ws.Range("I134").Value = 9409.5
$ws.Range("J134").Value = 3077.7693
$ws.Range("K134").Value = 28228.5
$ws.Range("L134").Value = 9233.3079
$ws.Range("M134").Value = -25693.5
$ws.Range("N134").Value = -14303.3079

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1436.5454
$ws.Range("I16").Value = 1378
$ws.Range("J16").Value = 1485.3334
$ws.Range("K16").Value = 1378
$ws.Range("L16").Value = 1485.3334
$ws.Range("M16").Value = -1091
$ws.Range("N16").Value = -2059.3334
# Row 99
$ws.Range("H99").Value = 10429220
$ws.Range("I99").Value = 12526.111
$ws.Range("J99").Value = 41679300
$ws.Range("K99").Value = 12526.111
$ws.Range("L99").Value = 41679300
$ws.Range("M99").Value = -11028.111
$ws.Range("N99").Value = -41682296
# Row 113
$ws.Range("H113").Value = 1436.5454
$ws.Range("I113").Value = 1378
$ws.Range("J113").Value = 1485.3334
$ws.Range("K113").Value = 1378
$ws.Range("L113").Value = 1485.3334
$ws.Range("M113").Value = 792
$ws.Range("N113").Value = -5825.3334
# Row 126
$ws.Range("H126").Value = 10429220
$ws.Range("I126").Value = 12526.111
$ws.Range("J126").Value = 41679300
$ws.Range("K126").Value = 37578.333
$ws.Range("L126").Value = 125037900
$ws.Range("M126").Value = -35108.333
$ws.Range("N126").Value = -125042840
# Row 134
$ws.Range("H134").Value = 2244.6597
$ws.Range("I134").Value = 2942.9644
$ws.Range("J134").Value = 1215.579
$ws.Range("K134").Value = 8828.893199999999
$ws.Range("L134").Value = 3646.737
$ws.Range("M134").Value = -6293.893199999999
$ws.Range("N134").Value = -8716.737000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Range("H103").Value = 2299.875
$ws.Range("J103").Value = 2774.75
$ws.Range("L103").Value = 8324.25
$ws.Range("N103").Value = -10082.25
# Row 113
$ws.Range("H113").Value = 1177004.9
$ws.Range("I113").Value = 1667200.4
$ws.Range("J113").Value = 476725.56
$ws.Range("K113").Value = 5001601.199999999
$ws.Range("L113").Value = 1430176.68
$ws.Range("M113").Value = -4999431.199999999
$ws.Range("N113").Value = -1434516.68

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 52632844
$ws.Range("I113").Value = 166667800
$ws.Range("J113").Value = 1328.2307
$ws.Range("K113").Value = 166667800
$ws.Range("L113").Value = 1328.2307
$ws.Range("M113").Value = -166665630
$ws.Range("N113").Value = -5668.2307
# Row 122
$ws.Range("H122").Value = 36719570
$ws.Range("I122").Value = 62637616
$ws.Range("J122").Value = 2335.25
$ws.Range("K122").Value = 187912848
$ws.Range("L122").Value = 7005.75
$ws.Range("M122").Value = -187910398
$ws.Range("N122").Value = -11905.75
# Row 132
$ws.Range("H132").Value = 2835.6775
$ws.Range("I132").Value = 2738.5
$ws.Range("J132").Value = 2897.0527
$ws.Range("K132").Value = 8215.5
$ws.Range("L132").Value = 8691.158100000001
$ws.Range("M132").Value = -5685.5
$ws.Range("N132").Value = -13751.1581

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1399.4814
$ws.Range("I61").Value = 1467.25
$ws.Range("J61").Value = 1300.909
$ws.Range("K61").Value = 1467.25
$ws.Range("L61").Value = 1300.909
$ws.Range("M61").Value = -1265.25
$ws.Range("N61").Value = -1704.909
# Row 113
$ws.Range("H113").Value = 1399.4814
$ws.Range("I113").Value = 1467.25
$ws.Range("J113").Value = 1300.909
$ws.Range("K113").Value = 1467.25
$ws.Range("L113").Value = 1300.909
$ws.Range("M113").Value = 702.75
$ws.Range("N113").Value = -5640.909
# Row 132
$ws.Range("H132").Value = 14847903
$ws.Range("I132").Value = 18430714
$ws.Range("J132").Value = 4829
$ws.Range("K132").Value = 55292142
$ws.Range("L132").Value = 14487
$ws.Range("M132").Value = -55289612
$ws.Range("N132").Value = -19547
# Row 136
$ws.Range("H136").Value = 9296.315000000001
$ws.Range("I136").Value = 6284
$ws.Range("J136").Value = 17730.8
$ws.Range("K136").Value = 18852
$ws.Range("L136").Value = 53192.39999999999
$ws.Range("M136").Value = -16302
$ws.Range("N136").Value = -58292.39999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4240
$ws.Range("I62").Value = 4240
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4240
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = -3616
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 4240
$ws.Range("I65").Value = 4240
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 21200
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = -18080
$ws.Range("M65").ClearContents()
# Row 107
$ws.Range("H107").Value = 115557336
$ws.Range("I107").Value = 200002600
$ws.Range("J107").Value = 10000751
$ws.Range("K107").Value = 600007800
$ws.Range("L107").Value = 30002253
$ws.Range("M107").Value = -600005880
$ws.Range("N107").Value = -30006093
# Row 113
$ws.Range("H113").Value = 1341.5
$ws.Range("I113").Value = 1391.6666
$ws.Range("J113").Value = 1281.3
$ws.Range("K113").Value = 4174.9998
$ws.Range("L113").Value = 3843.9
$ws.Range("M113").Value = -2004.9998
$ws.Range("N113").Value = -8183.9
# Row 122
$ws.Range("H122").Value = 2449.9062
$ws.Range("I122").Value = 2452.0435
$ws.Range("J122").Value = 2444.4443
$ws.Range("K122").Value = 7356.130500000001
$ws.Range("L122").Value = 7333.3329
$ws.Range("M122").Value = -4906.130500000001
$ws.Range("N122").Value = -12233.3329
# Row 132
$ws.Range("H132").Value = 21465.857
$ws.Range("I132").Value = 24558.191
$ws.Range("J132").Value = 2911.8572
$ws.Range("K132").Value = 73674.573
$ws.Range("L132").Value = 8735.571599999999
$ws.Range("M132").Value = -71144.573
$ws.Range("N132").Value = -13795.5716
# Row 136
$ws.Range("H136").Value = 7814964
$ws.Range("I136").Value = 2599.0732
$ws.Range("J136").Value = 21741354
$ws.Range("K136").Value = 7797.219599999999
$ws.Range("L136").Value = 65224062
$ws.Range("M136").Value = -5247.219599999999
$ws.Range("N136").Value = -65229162
